# Master Template files are added
#
# The shipped "User Shift" sheet had three example data rows (A3:A5 =
# vaishnavi / manisha / aditi) under the "UserId" header. This turns it
# into a blank, reusable master template: the sample rows are cleared
# (their cell formatting/style stays as-is), the leftover selection is
# reset back onto the sheet's live data cell, and the header row height
# is refreshed to match the cleaned-up template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Best-effort cosmetic touch-ups (safe no-ops if unsupported) --------
# Builtin "Hyperlink" cell style was relabeled to "Link" in the refreshed
# template's style table.
try {
    $wb.Styles.Item("Hyperlink").Name = "Link"
} catch {
}

# The app window was resized (recorded in the workbook's stored view).
try {
    $excel.ActiveWindow.Width = 19200
    $excel.ActiveWindow.Height = 6800
} catch {
}

# --- Core content edit: blank out the sample rows -----------------------
# Clear the three example user rows so the sheet is an empty template;
# ClearContents leaves the existing cell style/formatting (s="3") intact.
$ws.Range("A3:A5").ClearContents()

# The saved cursor position pointed at the now out-of-range A8; move it
# back onto A2 (the template's first real data cell under the header).
$ws.Range("A2").Select()

# Refresh the header row's height to match the template's row metrics.
$ws.Rows.Item(1).RowHeight = 15.5
